$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, pushing existing rows 53:79 down to 54:80
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new data record
$ws.Cells.Item(53, 1).Value2 = 3
$ws.Cells.Item(53, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(53, 3).Value2 = "Coquimbo"
$ws.Cells.Item(53, 4).Value2 = 44917
$ws.Cells.Item(53, 5).Value2 = 5
$ws.Cells.Item(53, 6).Value2 = 100112022
$ws.Cells.Item(53, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(53, 8).Value2 = "Perfection"
$ws.Cells.Item(53, 9).Value2 = "Primera"
$ws.Cells.Item(53, 10).Value2 = 38
$ws.Cells.Item(53, 11).Value2 = 28000
$ws.Cells.Item(53, 12).Value2 = 28000
$ws.Cells.Item(53, 13).Value2 = 28000
$ws.Cells.Item(53, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(53, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value2 = 1120
$ws.Cells.Item(53, 17).Value2 = 25
$ws.Cells.Item(53, 18).Value2 = "Hortaliza"
